$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ParticipantsTab ---
$ws.Range("A2").Value = "ParticipantsTab"
$participantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@
$ws.Range("B2").Value = $participantsQuery
$statQuery = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = "TC03_CDS_Filter_Study-Molecular Char Init_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC03_CDS_Filter_Study-Molecular Char Init_WebData.xlsx"

# --- Row 3: SamplesTab ---
$ws.Range("A3").Value = "SamplesTab"
$samplesQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
  ORDER By samp.sample_id 
  LIMIT 100
'@
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = "TC03_CDS_Filter_Study-Molecular Char Init_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC03_CDS_Filter_Study-Molecular Char Init_WebData.xlsx"

# --- Row 4: FilesTab ---
$ws.Range("A4").Value = "FilesTab"
$filesQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["MCI: Molecular Characterization Initiative"]
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name
Limit 100
'@
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = "TC03_CDS_Filter_Study-Molecular Char Init_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC03_CDS_Filter_Study-Molecular Char Init_WebData.xlsx"

# --- Row heights (content grew to the 409.5pt cap on rows 2 and 4) ---
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# --- Column A width: no longer best-fit, fixed width 19 ---
$ws.Columns.Item(1).ColumnWidth = 18.14

# --- Selection moves from A2 to B2 ---
$ws.Range("B2").Select()

Write-Host "Edit applied"
